$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1320.5385
$ws.Range("I137").Value = 1061.0454
$ws.Range("J137").Value = 1510.8334
$ws.Range("K137").Value = 3183.1362
$ws.Range("L137").Value = 4532.5002
$ws.Range("M137").Value = -633.1361999999999
$ws.Range("N137").Value = -9632.5002
$ws.Range("H138").Value = 3172.8818
$ws.Range("I138").Value = 1528.8286
$ws.Range("J138").Value = 4164.983
$ws.Range("K138").Value = 4586.4858
$ws.Range("L138").Value = 12494.949
$ws.Range("M138").Value = 553.5141999999996
$ws.Range("N138").Value = -22774.949

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11641.556
$ws.Range("I32").Value = 10617.058
$ws.Range("J32").Value = 47499
$ws.Range("K32").Value = 10617.058
$ws.Range("L32").Value = 47499
$ws.Range("M32").Value = -10330.058
$ws.Range("N32").Value = -48073
$ws.Range("H61").Value = 1685.7142
$ws.Range("I61").Value = 1583.5
$ws.Range("J61").Value = 1791.2258
$ws.Range("K61").Value = 1583.5
$ws.Range("L61").Value = 1791.2258
$ws.Range("M61").Value = -1371.5
$ws.Range("N61").Value = -2215.2258
$ws.Range("H74").Value = 1530.1578
$ws.Range("I74").Value = 1679.6904
$ws.Range("K74").Value = 1679.6904
$ws.Range("M74").Value = -805.6904
$ws.Range("H77").Value = 1530.1578
$ws.Range("I77").Value = 1679.6904
$ws.Range("K77").Value = 8398.451999999999
$ws.Range("M77").Value = -4030.451999999999
$ws.Range("H132").Value = 3270.4695
$ws.Range("I132").Value = 1001.37933
$ws.Range("J132").Value = 6560.65
$ws.Range("K132").Value = 3004.13799
$ws.Range("L132").Value = 19681.95
$ws.Range("M132").Value = -474.1379900000002
$ws.Range("N132").Value = -24741.95
$ws.Range("H136").Value = 1685.7142
$ws.Range("I136").Value = 1583.5
$ws.Range("J136").Value = 1791.2258
$ws.Range("K136").Value = 4750.5
$ws.Range("L136").Value = 5373.6774
$ws.Range("M136").Value = -2200.5
$ws.Range("N136").Value = -10473.6774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 552.6061
$ws.Range("I94").Value = 449.3913
$ws.Range("J94").Value = 790
$ws.Range("K94").Value = 449.3913
$ws.Range("L94").Value = 790
$ws.Range("M94").Value = 1.608699999999999
$ws.Range("N94").Value = -1692

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2582.4119
$ws.Range("I16").Value = 1233.4445
$ws.Range("J16").Value = 4100
$ws.Range("K16").Value = 1233.4445
$ws.Range("L16").Value = 4100
$ws.Range("M16").Value = -946.4445000000001
$ws.Range("N16").Value = -4674
$ws.Range("H31").Value = 2828.0422
$ws.Range("I31").Value = 1758.5
$ws.Range("J31").Value = 4059.6365
$ws.Range("K31").Value = 1758.5
$ws.Range("L31").Value = 4059.6365
$ws.Range("M31").Value = -1463.5
$ws.Range("N31").Value = -4649.636500000001
$ws.Range("H34").Value = 2828.0422
$ws.Range("I34").Value = 1758.5
$ws.Range("J34").Value = 4059.6365
$ws.Range("K34").Value = 1758.5
$ws.Range("L34").Value = 4059.6365
$ws.Range("M34").Value = -1556.5
$ws.Range("N34").Value = -4463.636500000001
$ws.Range("H113").Value = 2582.4119
$ws.Range("I113").Value = 1233.4445
$ws.Range("J113").Value = 4100
$ws.Range("K113").Value = 1233.4445
$ws.Range("L113").Value = 4100
$ws.Range("M113").Value = 936.5554999999999
$ws.Range("N113").Value = -8440
$ws.Range("H132").Value = 6805969.5
$ws.Range("I132").Value = 766.1389
$ws.Range("J132").Value = 25651148
$ws.Range("K132").Value = 2298.4167
$ws.Range("L132").Value = 76953444
$ws.Range("M132").Value = 231.5832999999998
$ws.Range("N132").Value = -76958504
$ws.Range("H134").Value = 1185.9824
$ws.Range("I134").Value = 1286.5227
$ws.Range("K134").Value = 3859.5681
$ws.Range("M134").Value = -1324.5681

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 8730707
$ws.Range("I113").Value = 11905233
$ws.Range("J113").Value = 7143445
$ws.Range("K113").Value = 35715699
$ws.Range("L113").Value = 21430335
$ws.Range("M113").Value = -35713529
$ws.Range("N113").Value = -21434675
$ws.Range("H131").Value = 733.79
$ws.Range("I131").Value = 410
$ws.Range("J131").Value = 769.76666
$ws.Range("K131").Value = 1230
$ws.Range("L131").Value = 2309.29998
$ws.Range("M131").Value = 3810
$ws.Range("N131").Value = -12389.29998
$ws.Range("H132").Value = 15629337
$ws.Range("J132").Value = 27784744
$ws.Range("L132").Value = 250062696
$ws.Range("N132").Value = -250067756

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5649.2964
$ws.Range("I132").Value = 1390.1666
$ws.Range("J132").Value = 9056.6
$ws.Range("K132").Value = 4170.4998
$ws.Range("L132").Value = 27169.8
$ws.Range("M132").Value = -1640.4998
$ws.Range("N132").Value = -32229.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2238.5625
$ws.Range("I7").Value = 2226
$ws.Range("K7").Value = 2226
$ws.Range("M7").Value = -2114
$ws.Range("H126").Value = 2238.5625
$ws.Range("I126").Value = 2226
$ws.Range("K126").Value = 6678
$ws.Range("M126").Value = -4208
$ws.Range("H132").Value = 5238.1294
$ws.Range("I132").Value = 1834.6897
$ws.Range("K132").Value = 5504.0691
$ws.Range("M132").Value = -2974.0691

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 50000864
$ws.Range("I81").Value = 62500580
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 125001160
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -125000099
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 50000864
$ws.Range("I84").Value = 62500580
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 625005800
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -625000496
$ws.Range("N84").Value = -30608
$ws.Range("H132").Value = 914.13
$ws.Range("I132").Value = 893.86304
$ws.Range("J132").Value = 968.9259
$ws.Range("K132").Value = 2681.58912
$ws.Range("L132").Value = 2906.7777
$ws.Range("M132").Value = -151.5891199999996
$ws.Range("N132").Value = -7966.7777
$ws.Range("H136").Value = 3492.1633
$ws.Range("I136").Value = 4458.6
$ws.Range("J136").Value = 1966.2106
$ws.Range("K136").Value = 13375.8
$ws.Range("L136").Value = 5898.6318
$ws.Range("M136").Value = -10825.8
$ws.Range("N136").Value = -10998.6318
